$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.46881016851428
$ws.Range("C2").Value = 13.54377188033832
$ws.Range("E2").Value = 16.77355237588091
$ws.Range("F2").Value = 47.50616550936569
$ws.Range("G2").Value = 3.669915681311299
$ws.Range("J2").Value = 9.265219716441752
$ws.Range("N2").Value = 18.36849905770418
$ws.Range("B3").Value = 19.84711276551646
$ws.Range("C3").Value = 12.94045050999947
$ws.Range("E3").Value = 16.72216121709578
$ws.Range("F3").Value = 47.07946520414222
$ws.Range("G3").Value = 3.674744770130402
$ws.Range("J3").Value = 9.282410435452039
$ws.Range("N3").Value = 18.44049805074787
$ws.Range("B4").Value = 19.46258293652042
$ws.Range("C4").Value = 12.56002308884923
$ws.Range("E4").Value = 16.69383337306526
$ws.Range("F4").Value = 46.83095792201943
$ws.Range("G4").Value = 3.677857012321019
$ws.Range("J4").Value = 9.294835203925411
$ws.Range("N4").Value = 18.4867694281206
$ws.Range("B5").Value = 19.30545007344866
$ws.Range("C5").Value = 12.40273950786562
$ws.Range("E5").Value = 16.68310941685904
$ws.Range("F5").Value = 46.73316768286469
$ws.Range("G5").Value = 3.679162449258026
$ws.Range("J5").Value = 9.300366507579843
$ws.Range("N5").Value = 18.50614470584586
$ws.Range("B6").Value = 19.27934030616031
$ws.Range("C6").Value = 12.37649448727461
$ws.Range("E6").Value = 16.68137846409373
$ws.Range("F6").Value = 46.71714215943163
$ws.Range("G6").Value = 3.679381466138781
$ws.Range("J6").Value = 9.301313188750406
$ws.Range("N6").Value = 18.50939333308645
$ws.Range("B7").Value = 19.46046517447551
$ws.Range("C7").Value = 12.55791068306315
$ws.Range("E7").Value = 16.69368541570903
$ws.Range("F7").Value = 46.82962489553097
$ws.Range("G7").Value = 3.677874467163876
$ws.Range("J7").Value = 9.294907908360836
$ws.Range("N7").Value = 18.48702862664971
$ws.Range("B8").Value = 20.25518920513057
$ws.Range("C8").Value = 13.33796217049593
$ws.Range("E8").Value = 16.75516838676933
$ws.Range("F8").Value = 47.35629410536921
$ws.Range("G8").Value = 3.671550309023925
$ws.Range("J8").Value = 9.270757847378452
$ws.Range("N8").Value = 18.392896215464
$ws.Range("B9").Value = 21.78027814904159
$ws.Range("C9").Value = 14.77879538740316
$ws.Range("E9").Value = 16.90092922054962
$ws.Range("F9").Value = 48.49194504020979
$ws.Range("G9").Value = 3.660308643623182
$ws.Range("J9").Value = 9.238322269682117
$ws.Range("N9").Value = 18.22465888030183
$ws.Range("B10").Value = 22.86681436179169
$ws.Range("C10").Value = 15.77235749047685
$ws.Range("E10").Value = 17.02278451348796
$ws.Range("F10").Value = 49.38297609184271
$ws.Range("G10").Value = 3.652745779350755
$ws.Range("J10").Value = 9.223706862297465
$ws.Range("N10").Value = 18.11099862756352
$ws.Range("B11").Value = 23.35123603449664
$ws.Range("C11").Value = 16.20845747469844
$ws.Range("E11").Value = 17.08128458731357
$ws.Range("F11").Value = 49.79924099662983
$ws.Range("G11").Value = 3.649454132782236
$ws.Range("J11").Value = 9.219082643316943
$ws.Range("N11").Value = 18.06144577369372
$ws.Range("B12").Value = 23.53307500145706
$ws.Range("C12").Value = 16.37119525448572
$ws.Range("E12").Value = 17.10386516807839
$ws.Range("F12").Value = 49.95831885318202
$ws.Range("G12").Value = 3.648228879813738
$ws.Range("J12").Value = 9.217624468525996
$ws.Range("N12").Value = 18.04299064795187
$ws.Range("B13").Value = 23.49398667580306
$ws.Range("C13").Value = 16.33625539450572
$ws.Range("E13").Value = 17.09898322319521
$ws.Range("F13").Value = 49.92399622436461
$ws.Range("G13").Value = 3.6484918187504
$ws.Range("J13").Value = 9.217925457143481
$ws.Range("N13").Value = 18.04695152461337
$ws.Range("B14").Value = 23.36622911442163
$ws.Range("C14").Value = 16.22189481021639
$ws.Range("E14").Value = 17.08313380796214
$ws.Range("F14").Value = 49.81230014134366
$ws.Range("G14").Value = 3.649352906112687
$ws.Range("J14").Value = 9.218956797866193
$ws.Range("N14").Value = 18.05992125681693
$ws.Range("B15").Value = 23.28776025566287
$ws.Range("C15").Value = 16.15152934300312
$ws.Range("E15").Value = 17.07348090233242
$ws.Range("F15").Value = 49.74406782154921
$ws.Range("G15").Value = 3.64988310579815
$ws.Range("J15").Value = 9.219626722413963
$ws.Range("N15").Value = 18.06790589595394
$ws.Range("B16").Value = 22.83494186094521
$ws.Range("C16").Value = 15.74352708631773
$ws.Range("E16").Value = 17.01902200477021
$ws.Range("F16").Value = 49.3559815736737
$ws.Range("G16").Value = 3.652963874984171
$ws.Range("J16").Value = 9.224049958788635
$ws.Range("N16").Value = 18.11428032450707
$ws.Range("B17").Value = 22.55449248194995
$ws.Range("C17").Value = 15.48907238681407
$ws.Range("E17").Value = 16.98638967408193
$ws.Range("F17").Value = 49.12061511633474
$ws.Range("G17").Value = 3.65489180432322
$ws.Range("J17").Value = 9.227283285007507
$ws.Range("N17").Value = 18.14328066567955
$ws.Range("B18").Value = 22.39226821108092
$ws.Range("C18").Value = 15.34122831241113
$ws.Range("E18").Value = 16.96790997896749
$ws.Range("F18").Value = 48.98627502364506
$ws.Range("G18").Value = 3.656014708198061
$ws.Range("J18").Value = 9.22933345249761
$ws.Range("N18").Value = 18.16016342320556
$ws.Range("B19").Value = 22.33719026398576
$ws.Range("C19").Value = 15.29091920590237
$ws.Range("E19").Value = 16.96170319678381
$ws.Range("F19").Value = 48.94097153784056
$ws.Range("G19").Value = 3.656397315429645
$ws.Range("J19").Value = 9.230060256265149
$ws.Range("N19").Value = 18.16591441878451
$ws.Range("B20").Value = 22.58444312165787
$ws.Range("C20").Value = 15.51631451046699
$ws.Range("E20").Value = 16.98983356711713
$ws.Range("F20").Value = 49.14556390009374
$ws.Range("G20").Value = 3.654685124185644
$ws.Range("J20").Value = 9.226919368317946
$ws.Range("N20").Value = 18.140172568823
$ws.Range("B21").Value = 23.40379943809984
$ws.Range("C21").Value = 16.25555139336961
$ws.Range("E21").Value = 17.08777766144854
$ws.Range("F21").Value = 49.84506971065188
$ws.Range("G21").Value = 3.649099409244018
$ws.Range("J21").Value = 9.21864590444248
$ws.Range("N21").Value = 18.05610333265616
$ws.Range("B22").Value = 23.929891518751
$ws.Range("C22").Value = 16.72463105815748
$ws.Range("E22").Value = 17.15427721320884
$ws.Range("F22").Value = 50.31061438432866
$ws.Range("G22").Value = 3.645572443222911
$ws.Range("J22").Value = 9.214946835566472
$ws.Range("N22").Value = 18.00296308485878
$ws.Range("B23").Value = 23.65002306413638
$ws.Range("C23").Value = 16.47559548291947
$ws.Range("E23").Value = 17.1185620725214
$ws.Range("F23").Value = 50.06141959316206
$ws.Range("G23").Value = 3.647443594968449
$ws.Range("J23").Value = 9.21676423054825
$ws.Range("N23").Value = 18.03115995256091
$ws.Range("B24").Value = 22.57090550744092
$ws.Range("C24").Value = 15.50400317947299
$ws.Range("E24").Value = 16.98827570662097
$ws.Range("F24").Value = 49.13428150796085
$ws.Range("G24").Value = 3.654778519062893
$ws.Range("J24").Value = 9.227083299380768
$ws.Range("N24").Value = 18.14157708507758
$ws.Range("B25").Value = 21.37275023028456
$ws.Range("C25").Value = 14.39974788891367
$ws.Range("E25").Value = 16.85885833503456
$ws.Range("F25").Value = 48.17431057901867
$ws.Range("G25").Value = 3.663226740433583
$ws.Range("J25").Value = 9.245486917442845
$ws.Range("N25").Value = 18.2684235234349
